$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.092628121376038
$ws.Range("B1").Value = 1.962450504302979
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.0750732421875
$ws.Range("E1").Value = 1.13274621963501
